$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Bump every value in D2:D23 by 90 (e.g. -50 -> 40, -60 -> 30, -48 -> 42)
for ($r = 2; $r -le 23; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.Value2 = $cell.Value2 + 90
}

# Update the active selection on the survey sheet to D2:D23 (active cell D2)
$ws.Activate()
$ws.Range("D2:D23").Select()
